# Sprint 4 burn-up sheet update: fill in the "Dag 9" (N) and "Dag 10" (O)
# columns for each backlog item (rows 2-21) with the day's remaining
# hours. Most rows simply repeat the "Dag 8" (M) value; a few items
# (rows 8, 13, 18) got extra hours logged on day 10 (column O).
#
# Rows 10-18/20/21 had never been touched before, so their N/O cells
# were still sitting on a stray duplicate style; fix that up to match
# the rest of the row (same style as column M) before writing values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# row -> [N value, O value]
$dayValues = [ordered]@{
  2  = @(5.5, 5.5)
  3  = @(2,   2)
  4  = @(2,   2)
  5  = @(2.5, 2.5)
  6  = @(2.5, 2.5)
  7  = @(2,   2)
  8  = @(0,   2.5)
  9  = @(1,   1)
  10 = @(3,   3)
  11 = @(3,   3)
  12 = @(0,   0)
  13 = @(0,   2)
  14 = @(0,   0)
  15 = @(0,   0)
  16 = @(5,   5)
  17 = @(0,   0)
  18 = @(5.5, 10)
  19 = @(3,   3)
  20 = @(16,  16)
  21 = @(13,  13)
}

# Rows whose N/O cells are still on the old, unused duplicate style and
# need reformatting to match column M before the values go in.
$rowsNeedingRestyle = @(10, 11, 12, 13, 14, 15, 16, 17, 18, 20, 21)

foreach ($row in $dayValues.Keys) {
  $pair = $dayValues[$row]

  if ($rowsNeedingRestyle -contains $row) {
    $ws.Range("M$row").Copy()
    $ws.Range("N${row}:O${row}").PasteSpecial(-4122)
  }

  $ws.Range("N$row").Value = $pair[0]
  $ws.Range("O$row").Value = $pair[1]
}

# Extend the "Sum" row's running totals into the newly-populated columns.
$ws.Range("M23").Formula = "=SUM(M1:M21)"
$ws.Range("N23").Formula = "=SUM(N1:N21)"
$ws.Range("O23").Formula = "=SUM(O1:O21)"

# Leave the selection where the editor (SM) ended up.
$ws.Range("Q17").Select()
